$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 800.2143
$ws.Cells.Item(5, 9).Value = 1487.2858
$ws.Cells.Item(5, 10).Value = 113.14286
$ws.Cells.Item(5, 11).Value = 1487.2858
$ws.Cells.Item(5, 12).Value = 113.14286
$ws.Cells.Item(5, 13).Value = -1372.2858
$ws.Cells.Item(5, 14).Value = -343.14286
$ws.Cells.Item(13, 8).Value = 205
$ws.Cells.Item(13, 9).Value = 205
$ws.Cells.Item(13, 11).Value = 205
$ws.Cells.Item(13, 13).Value = -36
$ws.Cells.Item(55, 8).Value = 1043.3334
$ws.Cells.Item(55, 9).Value = 991.75
$ws.Cells.Item(55, 10).Value = 1112.1111
$ws.Cells.Item(55, 11).Value = 991.75
$ws.Cells.Item(55, 12).Value = 1112.1111
$ws.Cells.Item(55, 13).Value = -777.75
$ws.Cells.Item(55, 14).Value = -1540.1111
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 14).ClearContents()
$ws.Cells.Item(64, 8).Value = 65312.25
$ws.Cells.Item(64, 9).Value = 127574.75
$ws.Cells.Item(64, 10).Value = 3049.75
$ws.Cells.Item(64, 11).Value = 127574.75
$ws.Cells.Item(64, 12).Value = 3049.75
$ws.Cells.Item(64, 13).Value = -127326.75
$ws.Cells.Item(64, 14).Value = -3545.75
$ws.Cells.Item(67, 8).Value = 65312.25
$ws.Cells.Item(67, 9).Value = 127574.75
$ws.Cells.Item(67, 10).Value = 3049.75
$ws.Cells.Item(67, 11).Value = 127574.75
$ws.Cells.Item(67, 12).Value = 3049.75
$ws.Cells.Item(67, 13).Value = -126716.75
$ws.Cells.Item(67, 14).Value = -4765.75
$ws.Cells.Item(129, 8).Value = 635573.8
$ws.Cells.Item(129, 9).Value = 2002707.9
$ws.Cells.Item(129, 10).Value = 4588.846
$ws.Cells.Item(129, 11).Value = 6008123.699999999
$ws.Cells.Item(129, 12).Value = 13766.538
$ws.Cells.Item(129, 13).Value = -6003123.699999999
$ws.Cells.Item(129, 14).Value = -23766.538
$ws.Cells.Item(132, 8).Value = 23262.25
$ws.Cells.Item(132, 9).Value = 3783.2104
$ws.Cells.Item(132, 10).Value = 146629.5
$ws.Cells.Item(132, 11).Value = 11349.6312
$ws.Cells.Item(132, 12).Value = 439888.5
$ws.Cells.Item(132, 13).Value = -8819.6312
$ws.Cells.Item(132, 14).Value = -444948.5
$ws.Cells.Item(138, 8).Value = 2645.5789
$ws.Cells.Item(138, 9).Value = 2032.4062
$ws.Cells.Item(138, 10).Value = 3091.5227
$ws.Cells.Item(138, 11).Value = 6097.2186
$ws.Cells.Item(138, 12).Value = 9274.5681
$ws.Cells.Item(138, 13).Value = -957.2186000000002
$ws.Cells.Item(138, 14).Value = -19554.5681
$ws.Cells.Item(139, 8).Value = 49712.06
$ws.Cells.Item(139, 10).Value = 49712.06
$ws.Cells.Item(139, 12).Value = 49712.06
$ws.Cells.Item(139, 14).Value = -59992.06
$ws.Cells.Item(140, 8).Value = 57847.145
$ws.Cells.Item(140, 10).Value = 57847.145
$ws.Cells.Item(140, 12).Value = 57847.145
$ws.Cells.Item(140, 14).Value = -68207.14499999999
$ws.Cells.Item(141, 8).Value = 2607.8333
$ws.Cells.Item(141, 9).Value = 2051.8096
$ws.Cells.Item(141, 10).Value = 6500
$ws.Cells.Item(141, 11).Value = 6155.4288
$ws.Cells.Item(141, 12).Value = 19500
$ws.Cells.Item(141, 13).Value = -975.4287999999997
$ws.Cells.Item(141, 14).Value = -29860

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(135, 8).Value = 31697.5
$ws.Cells.Item(135, 10).Value = 31697.5
$ws.Cells.Item(135, 12).Value = 31697.5
$ws.Cells.Item(135, 14).Value = -41837.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 300
$ws.Cells.Item(2, 9).Value = 500
$ws.Cells.Item(2, 10).Value = 100
$ws.Cells.Item(2, 11).Value = 500
$ws.Cells.Item(2, 12).Value = 100
$ws.Cells.Item(2, 13).Value = -387
$ws.Cells.Item(2, 14).Value = -326
$ws.Cells.Item(6, 8).Value = 2290557
$ws.Cells.Item(6, 9).Value = 6668667
$ws.Cells.Item(6, 10).Value = 101502
$ws.Cells.Item(6, 11).Value = 6668667
$ws.Cells.Item(6, 12).Value = 101502
$ws.Cells.Item(6, 13).Value = -6668554
$ws.Cells.Item(6, 14).Value = -101728
$ws.Cells.Item(7, 8).Value = 600
$ws.Cells.Item(7, 9).Value = 433.33334
$ws.Cells.Item(7, 11).Value = 433.33334
$ws.Cells.Item(7, 13).Value = -320.33334
$ws.Cells.Item(11, 8).Value = 12795
$ws.Cells.Item(11, 9).Value = 1005
$ws.Cells.Item(11, 10).Value = 16725
$ws.Cells.Item(11, 11).Value = 1005
$ws.Cells.Item(11, 12).Value = 16725
$ws.Cells.Item(11, 13).Value = -865
$ws.Cells.Item(11, 14).Value = -17005
$ws.Cells.Item(20, 8).Value = 49816.75
$ws.Cells.Item(20, 10).Value = 49816.75
$ws.Cells.Item(20, 12).Value = 49816.75
$ws.Cells.Item(20, 14).Value = -50288.75
$ws.Cells.Item(30, 8).Value = 49816.75
$ws.Cells.Item(30, 10).Value = 49816.75
$ws.Cells.Item(30, 12).Value = 49816.75
$ws.Cells.Item(30, 14).Value = -49998.75
$ws.Cells.Item(62, 8).Value = 2679.9565
$ws.Cells.Item(62, 9).Value = 2541.0557
$ws.Cells.Item(62, 11).Value = 2541.0557
$ws.Cells.Item(62, 13).Value = -1917.0557
$ws.Cells.Item(65, 8).Value = 2679.9565
$ws.Cells.Item(65, 9).Value = 2541.0557
$ws.Cells.Item(65, 11).Value = 12705.2785
$ws.Cells.Item(65, 13).Value = -9585.2785
$ws.Cells.Item(128, 8).Value = 49816.75
$ws.Cells.Item(128, 10).Value = 49816.75
$ws.Cells.Item(128, 12).Value = 49816.75
$ws.Cells.Item(128, 14).Value = -59776.75
$ws.Cells.Item(130, 8).Value = 30687
$ws.Cells.Item(130, 10).Value = 30687
$ws.Cells.Item(130, 12).Value = 30687
$ws.Cells.Item(130, 14).Value = -40727
$ws.Cells.Item(133, 8).Value = 35636.617
$ws.Cells.Item(133, 10).Value = 35636.617
$ws.Cells.Item(133, 12).Value = 35636.617
$ws.Cells.Item(133, 14).Value = -40696.617
$ws.Cells.Item(134, 8).Value = 3233.5356
$ws.Cells.Item(134, 9).Value = 1446.7222
$ws.Cells.Item(134, 10).Value = 6449.8
$ws.Cells.Item(134, 11).Value = 4340.1666
$ws.Cells.Item(134, 12).Value = 19349.4
$ws.Cells.Item(134, 13).Value = -1805.1666
$ws.Cells.Item(134, 14).Value = -24419.4
$ws.Cells.Item(135, 8).Value = 53965.08
$ws.Cells.Item(135, 10).Value = 53965.08
$ws.Cells.Item(135, 12).Value = 53965.08
$ws.Cells.Item(135, 14).Value = -64105.08
$ws.Cells.Item(137, 8).Value = 28543.4
$ws.Cells.Item(137, 10).Value = 28543.4
$ws.Cells.Item(137, 12).Value = 28543.4
$ws.Cells.Item(137, 14).Value = -38743.4
$ws.Cells.Item(138, 8).Value = 36530.066
$ws.Cells.Item(138, 10).Value = 36530.066
$ws.Cells.Item(138, 12).Value = 36530.066
$ws.Cells.Item(138, 14).Value = -46810.066
$ws.Cells.Item(139, 8).Value = 48963.3
$ws.Cells.Item(139, 10).Value = 49737
$ws.Cells.Item(139, 12).Value = 49737
$ws.Cells.Item(139, 14).Value = -60017
$ws.Cells.Item(140, 8).Value = 55704.453
$ws.Cells.Item(140, 10).Value = 55704.453
$ws.Cells.Item(140, 12).Value = 55704.453
$ws.Cells.Item(140, 14).Value = -66064.45300000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 9563855
$ws.Cells.Item(2, 9).Value = 61
$ws.Cells.Item(2, 10).Value = 15189616
$ws.Cells.Item(2, 11).Value = 366
$ws.Cells.Item(2, 12).Value = 91137696
$ws.Cells.Item(2, 13).Value = -253
$ws.Cells.Item(2, 14).Value = -91137922
$ws.Cells.Item(131, 8).Value = 4830.069
$ws.Cells.Item(131, 9).Value = 25599.75
$ws.Cells.Item(131, 10).Value = 1506.92
$ws.Cells.Item(131, 11).Value = 76799.25
$ws.Cells.Item(131, 12).Value = 4520.76
$ws.Cells.Item(131, 13).Value = -71759.25
$ws.Cells.Item(131, 14).Value = -14600.76

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2388.8
$ws.Cells.Item(68, 9).Value = 2294.3333
$ws.Cells.Item(68, 10).Value = 2766.6667
$ws.Cells.Item(68, 11).Value = 2294.3333
$ws.Cells.Item(68, 12).Value = 2766.6667
$ws.Cells.Item(68, 13).Value = -1545.3333
$ws.Cells.Item(68, 14).Value = -4264.6667
$ws.Cells.Item(71, 8).Value = 2388.8
$ws.Cells.Item(71, 9).Value = 2294.3333
$ws.Cells.Item(71, 10).Value = 2766.6667
$ws.Cells.Item(71, 11).Value = 11471.6665
$ws.Cells.Item(71, 12).Value = 13833.3335
$ws.Cells.Item(71, 13).Value = -7727.666499999999
$ws.Cells.Item(71, 14).Value = -21321.3335

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 1200
$ws.Cells.Item(2, 10).Value = 1200
$ws.Cells.Item(2, 12).Value = 1200
$ws.Cells.Item(2, 14).Value = -1424
$ws.Cells.Item(4, 8).Value = 1950
$ws.Cells.Item(4, 10).Value = 1950
$ws.Cells.Item(4, 12).Value = 1950
$ws.Cells.Item(4, 14).Value = -2176
$ws.Cells.Item(6, 8).Value = 1390219.5
$ws.Cells.Item(6, 9).Value = 5555633
$ws.Cells.Item(6, 10).Value = 1748.3334
$ws.Cells.Item(6, 11).Value = 5555633
$ws.Cells.Item(6, 12).Value = 1748.3334
$ws.Cells.Item(6, 13).Value = -5555518
$ws.Cells.Item(6, 14).Value = -1978.3334
$ws.Cells.Item(7, 8).Value = 100005
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 100005
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 100005
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(7, 14).Value = -100231
$ws.Cells.Item(8, 8).Value = 74802.39999999999
$ws.Cells.Item(8, 10).Value = 74802.39999999999
$ws.Cells.Item(8, 12).Value = 74802.39999999999
$ws.Cells.Item(8, 14).Value = -75082.39999999999
$ws.Cells.Item(9, 8).Value = 61000.2
$ws.Cells.Item(9, 9).Value = 2000
$ws.Cells.Item(9, 10).Value = 75750.25
$ws.Cells.Item(9, 11).Value = 2000
$ws.Cells.Item(9, 12).Value = 75750.25
$ws.Cells.Item(9, 13).Value = -1860
$ws.Cells.Item(9, 14).Value = -76030.25
$ws.Cells.Item(10, 8).Value = 2995
$ws.Cells.Item(10, 10).Value = 2995
$ws.Cells.Item(10, 12).Value = 2995
$ws.Cells.Item(10, 14).Value = -3333
$ws.Cells.Item(13, 8).Value = 1119.75
$ws.Cells.Item(13, 9).Value = 999.5
$ws.Cells.Item(13, 10).Value = 1240
$ws.Cells.Item(13, 11).Value = 999.5
$ws.Cells.Item(13, 12).Value = 1240
$ws.Cells.Item(13, 13).Value = -859.5
$ws.Cells.Item(13, 14).Value = -1520
$ws.Cells.Item(62, 8).Value = 2444.6667
$ws.Cells.Item(62, 9).Value = 1751
$ws.Cells.Item(62, 10).Value = 2642.8572
$ws.Cells.Item(62, 11).Value = 1751
$ws.Cells.Item(62, 12).Value = 2642.8572
$ws.Cells.Item(62, 13).Value = -1127
$ws.Cells.Item(62, 14).Value = -3890.8572
$ws.Cells.Item(65, 8).Value = 2444.6667
$ws.Cells.Item(65, 9).Value = 1751
$ws.Cells.Item(65, 10).Value = 2642.8572
$ws.Cells.Item(65, 11).Value = 8755
$ws.Cells.Item(65, 12).Value = 13214.286
$ws.Cells.Item(65, 13).Value = -5635
$ws.Cells.Item(65, 14).Value = -19454.286
